# Refresh the cryptocurrency price/volume table with latest scraped values.
# (Equivalent to the "Updated cryptos list ... with GitHub Actions" commit.)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "28.679.17"
$ws.Range("E2").Value = "  +1.48%  "
# Row 3
$ws.Range("D3").Value = "1.573.29"
$ws.Range("E3").Value = "  -0.83%  "
# Row 4
$ws.Range("E4").Value = "  +0.22%  "
# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "213.36"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.02%  "
# Row 6
$ws.Range("E6").Value = "  +0.06%  "
# Row 7
$ws.Range("E7").Value = "  +0.29%  "
# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "44.54"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +1.35%  "
# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "24.09"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +0.48%  "
# Row 10
$ws.Range("E10").Value = "  -1.16%  "
# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0592"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -0.89%  "
# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.0890"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +0.53%  "
# Row 13
$ws.Range("D13").Value = "1.798.04"
$ws.Range("E13").Value = "  -0.84%  "
# Row 14
$ws.Range("D14").Value = "1.575.56"
$ws.Range("E14").Value = "  -0.67%  "
# Row 15
$ws.Range("B15").Value = "WrappedBTC"
$ws.Range("C15").Value = "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
$ws.Range("D15").Value = "28.674.49"
$ws.Range("E15").Value = "  +1.41%  "
# Row 16
$ws.Range("B16").Value = "Polygon"
$ws.Range("C16").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.521"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -1.63%  "
# Row 17
$ws.Range("E17").Value = "  -1.79%  "
# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "62.31"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -1.41%  "
# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "230.00"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +1.11%  "
# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "7.36"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -1.44%  "
# Row 21
$ws.Range("D21").Value = "0.0₃0691"
$ws.Range("E21").Value = "  -2.02%  "
# Row 22
$ws.Range("E22").Value = "  +0.10%  "
# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "3.88"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -4.65%  "
# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "9.16"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -1.66%  "
# Row 25
$ws.Range("E25").Value = "  +4.84%  "
# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "151.83"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -0.04%  "
# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "15.00"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -0.89%  "
# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "6.45"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -1.65%  "
# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.104"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -2.38%  "
# Row 30
$ws.Range("E30").Value = "  +0.25%  "
# Row 31
$ws.Range("E31").Value = "  +1.99%  "
# Row 32
$ws.Range("E32").Value = "  -2.21%  "
# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "3.21"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -0.80%  "
# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "3.11"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -1.33%  "
# Row 35
$ws.Range("D35").Value = "1.395.57"
$ws.Range("E35").Value = "  -0.27%  "
# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.04"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +2.24%  "
# Row 37
$ws.Range("E37").Value = "  -3.80%  "
# Row 38
$ws.Range("E38").Value = "  +0.77%  "
# Row 39
$ws.Range("E39").Value = "  +2.94%  "
# Row 40
$ws.Range("E40").Value = "  -0.46%  "
# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.523"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -3.04%  "
# Row 42
$ws.Range("E42").Value = "  +0.30%  "
# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.89"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +0.44%  "
# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.792"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -2.23%  "
# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0463"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -0.10%  "
# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "5.48"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -2.19%  "
# Row 47
$ws.Range("E47").Value = "  -1.79%  "
# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "63.08"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -1.78%  "
# Row 49
$ws.Range("D49").Value = "1.710.42"
$ws.Range("E49").Value = "  -0.68%  "
# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "86.51"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +0.04%  "
# Row 51
$ws.Range("D51").Value = "0.0₆0102"
$ws.Range("E51").Value = "  -0.65%  "
